$d = $word.ActiveDocument

# The title, author and abstract paragraphs currently have their text split
# across several runs (one run per word, with separate single-space runs in
# between). Collapse each of those paragraphs down to a single run holding
# the full sentence, leaving every other paragraph (including the "Guide: ..."
# hyperlink text elsewhere in the document) untouched.

function Set-ParagraphText($paragraph, [string]$text) {
    $r = $paragraph.Range
    # Exclude the trailing paragraph mark from the search/replace range.
    $r.MoveEnd(1, -1) | Out-Null
    $r.Find.ClearFormatting()
    $r.Find.Execute($text, $true, $false, $false, $false, $false, `
                     $true, 1, $false, $text, 2) | Out-Null
}

foreach ($p in $d.Paragraphs) {
    $styleName = $p.Range.Style.NameLocal
    switch ($styleName) {
        "Title" {
            Set-ParagraphText $p "Questions: Further sigma notation"
        }
        "Author" {
            Set-ParagraphText $p "Ifan Howells-Baines, Mark Toner"
        }
        "Abstract" {
            Set-ParagraphText $p "Questions relating to the guide on sigma notation"
        }
    }
}
